$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 25:50 hold the US-listed companies. Column A currently holds the
# Korean display name (the value we now want to preserve in column C),
# column B holds the English legal name (unchanged), and column C
# sometimes holds a short alias (AMD/P&G/J&J/머크/에이엠디/페이스북/구글)
# that is being retired in favour of the SEC CIK code now going into A.
$cikByRow = @{
    25 = "0000002488"
    26 = "0000019617"
    27 = "0000034088"
    28 = "0000059478"
    29 = "0000070858"
    30 = "0000080424"
    31 = "0000104169"
    32 = "0000200406"
    33 = "0000310158"
    34 = "0000320193"
    35 = "0000354950"
    36 = "0000731766"
    37 = "0000789019"
    38 = "0000909832"
    39 = "0001018724"
    40 = "0001045810"
    41 = "0001067983"
    42 = "0001108524"
    43 = "0001141391"
    44 = "0001318605"
    45 = "0001326801"
    46 = "0001341439"
    47 = "0001403161"
    48 = "0001551152"
    49 = "0001652044"
    50 = "0001730168"
}

# Make column A hold text (so the zero-padded CIK codes keep their
# leading zeros) before writing into it.
$ws.Range("A25:A50").NumberFormat = "@"

for ($r = 25; $r -le 50; $r++) {
    $korName = $ws.Cells.Item($r, 1).Value()
    # Move the old Korean display name from A into C (replacing whatever
    # alias used to live there, including the now-removed P&G/J&J/머크/
    # 에이엠디/페이스북/구글 entries).
    $ws.Cells.Item($r, 3).Value = $korName
    $ws.Cells.Item($r, 3).NumberFormat = "General"
    # Put the finalized SEC CIK code into A.
    $ws.Cells.Item($r, 1).Value = $cikByRow[$r]
}
